$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46075 -> 46076) for every data row (rows 2 through 120).
$range = $ws.Range("C2:C120")
$range.Value = 46076
